$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "67.050.39"
$c.ClearFormats()
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  -0.64%  "
$c.ClearFormats()
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.602.35"
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  +0.01%  "
$c.ClearFormats()
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  +0.02%  "
$c.ClearFormats()
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "591.37"
$c.ClearFormats()
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  -0.54%  "
$c.ClearFormats()
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "151.49"
$c.ClearFormats()
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  -1.88%  "
$c.ClearFormats()
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  -0.04%  "
$c.ClearFormats()
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  +0.82%  "
$c.ClearFormats()
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.600.08"
$c.ClearFormats()
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  +0.07%  "
$c.ClearFormats()
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.122"
$c.ClearFormats()
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  -3.13%  "
$c.ClearFormats()
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  +0.42%  "
$c.ClearFormats()
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "5.18"
$c.ClearFormats()
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  -1.59%  "
$c.ClearFormats()
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.346"
$c.ClearFormats()
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  -2.32%  "
$c.ClearFormats()
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "27.51"
$c.ClearFormats()
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  -0.31%  "
$c.ClearFormats()
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.073.69"
$c.ClearFormats()
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  -0.21%  "
$c.ClearFormats()
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  -4.09%  "
$c.ClearFormats()
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "66.909.94"
$c.ClearFormats()
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  -0.92%  "
$c.ClearFormats()
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.597.77"
$c.ClearFormats()
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  -0.24%  "
$c.ClearFormats()
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "362.93"
$c.ClearFormats()
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  -0.05%  "
$c.ClearFormats()
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "10.99"
$c.ClearFormats()
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  -1.30%  "
$c.ClearFormats()
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "7.35"
$c.ClearFormats()
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  -4.61%  "
$c.ClearFormats()
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  -0.38%  "
$c.ClearFormats()
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "2.06"
$c.ClearFormats()
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  +1.36%  "
$c.ClearFormats()
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  +0.10%  "
$c.ClearFormats()
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "10.07"
$c.ClearFormats()
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  +2.19%  "
$c.ClearFormats()
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "66.70"
$c.ClearFormats()
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  -0.85%  "
$c.ClearFormats()
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.729.28"
$c.ClearFormats()
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  -0.76%  "
$c.ClearFormats()
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "583.81"
$c.ClearFormats()
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  +0.96%  "
$c.ClearFormats()
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.ClearFormats()
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  -0.46%  "
$c.ClearFormats()
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  -3.35%  "
$c.ClearFormats()
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  -2.87%  "
$c.ClearFormats()
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  -3.43%  "
$c.ClearFormats()
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  -3.19%  "
$c.ClearFormats()
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  +0.07%  "
$c.ClearFormats()
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  -1.75%  "
$c.ClearFormats()
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  -1.96%  "
$c.ClearFormats()
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "156.30"
$c.ClearFormats()
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  +0.83%  "
$c.ClearFormats()
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "18.83"
$c.ClearFormats()
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -2.64%  "
$c.ClearFormats()
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  -0.98%  "
$c.ClearFormats()
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "5.22"
$c.ClearFormats()
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  -2.49%  "
$c.ClearFormats()
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  -3.52%  "
$c.ClearFormats()
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.57"
$c.ClearFormats()
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  +0.65%  "
$c.ClearFormats()
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "16.77"
$c.ClearFormats()
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  +1.63%  "
$c.ClearFormats()
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  -0.02%  "
$c.ClearFormats()
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0₆0295"
$c.ClearFormats()
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  +1.09%  "
$c.ClearFormats()
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "153.68"
$c.ClearFormats()
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  -1.53%  "
$c.ClearFormats()
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  -0.39%  "
$c.ClearFormats()
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "21.64"
$c.ClearFormats()
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  +3.68%  "
$c.ClearFormats()
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.615"
$c.ClearFormats()
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  -2.11%  "
$c.ClearFormats()
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.550"
$c.ClearFormats()
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  -2.67%  "
$c.ClearFormats()
